$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1270.7142
$ws.Range("J98").Value = 2000
$ws.Range("L98").Value = 2000
$ws.Range("N98").Value = -4996
$ws.Range("H113").Value = 2126.25
$ws.Range("I113").Value = 1951.6666
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 1951.6666
$ws.Range("L113").Value = 2650
$ws.Range("M113").Value = 1302.3334
$ws.Range("N113").Value = -9158
$ws.Range("H122").Value = 1270.7142
$ws.Range("J122").Value = 2000
$ws.Range("L122").Value = 6000
$ws.Range("N122").Value = -10900
$ws.Range("H138").Value = 2231.625
$ws.Range("J138").Value = 2094.5305
$ws.Range("L138").Value = 6283.5915
$ws.Range("N138").Value = -16563.5915

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1677.2354
$ws.Range("I2").Value = 1864.1538
$ws.Range("J2").Value = 1069.75
$ws.Range("K2").Value = 1864.1538
$ws.Range("L2").Value = 1069.75
$ws.Range("M2").Value = -1751.1538
$ws.Range("N2").Value = -1295.75
$ws.Range("H32").Value = 403906.22
$ws.Range("I32").Value = 471091.6
$ws.Range("J32").Value = 17590.334
$ws.Range("K32").Value = 471091.6
$ws.Range("L32").Value = 17590.334
$ws.Range("M32").Value = -470804.6
$ws.Range("N32").Value = -18164.334
$ws.Range("H63").Value = 4734.933
$ws.Range("I63").Value = 2943.5
$ws.Range("J63").Value = 6782.2856
$ws.Range("K63").Value = 2943.5
$ws.Range("L63").Value = 6782.2856
$ws.Range("M63").Value = -2257.5
$ws.Range("N63").Value = -8154.2856
$ws.Range("H66").Value = 4734.933
$ws.Range("I66").Value = 2943.5
$ws.Range("J66").Value = 6782.2856
$ws.Range("K66").Value = 14717.5
$ws.Range("L66").Value = 33911.428
$ws.Range("M66").Value = -11285.5
$ws.Range("N66").Value = -40775.428
$ws.Range("H74").Value = 1655.8636
$ws.Range("I74").Value = 1313.8334
$ws.Range("J74").Value = 3195
$ws.Range("K74").Value = 1313.8334
$ws.Range("L74").Value = 3195
$ws.Range("M74").Value = -439.8334
$ws.Range("N74").Value = -4943
$ws.Range("H77").Value = 1655.8636
$ws.Range("I77").Value = 1313.8334
$ws.Range("J77").Value = 3195
$ws.Range("K77").Value = 6569.166999999999
$ws.Range("L77").Value = 15975
$ws.Range("M77").Value = -2201.166999999999
$ws.Range("N77").Value = -24711
$ws.Range("H88").Value = 3105
$ws.Range("I88").Value = 3129.5
$ws.Range("K88").Value = 3129.5
$ws.Range("M88").Value = -2723.5
$ws.Range("H91").Value = 3105
$ws.Range("I91").Value = 3129.5
$ws.Range("K91").Value = 3129.5
$ws.Range("M91").Value = -1725.5
$ws.Range("H102").Value = 2437.5
$ws.Range("I102").Value = 2250
$ws.Range("K102").Value = 2250
$ws.Range("M102").Value = -628
$ws.Range("H116").Value = 1677.2354
$ws.Range("I116").Value = 1864.1538
$ws.Range("J116").Value = 1069.75
$ws.Range("K116").Value = 1864.1538
$ws.Range("L116").Value = 1069.75
$ws.Range("M116").Value = 429.8462
$ws.Range("N116").Value = -5657.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1677.2354
$ws.Range("I3").Value = 1864.1538
$ws.Range("J3").Value = 1069.75
$ws.Range("K3").Value = 1864.1538
$ws.Range("L3").Value = 1069.75
$ws.Range("M3").Value = -1750.1538
$ws.Range("N3").Value = -1297.75
$ws.Range("H86").Value = 55557670
$ws.Range("I86").Value = 66668812
$ws.Range("K86").Value = 66668812
$ws.Range("M86").Value = -66667689
$ws.Range("H89").Value = 55557670
$ws.Range("I89").Value = 66668812
$ws.Range("K89").Value = 333344060
$ws.Range("M89").Value = -333338444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 30500.25
$ws.Range("J4").Value = 30500.25
$ws.Range("L4").Value = 30500.25
$ws.Range("N4").Value = -30724.25
$ws.Range("H31").Value = 6057.4375
$ws.Range("I31").Value = 1519.1111
$ws.Range("J31").Value = 11892.429
$ws.Range("K31").Value = 1519.1111
$ws.Range("L31").Value = 11892.429
$ws.Range("M31").Value = -1224.1111
$ws.Range("N31").Value = -12482.429
$ws.Range("H34").Value = 6057.4375
$ws.Range("I34").Value = 1519.1111
$ws.Range("J34").Value = 11892.429
$ws.Range("K34").Value = 1519.1111
$ws.Range("L34").Value = 11892.429
$ws.Range("M34").Value = -1317.1111
$ws.Range("N34").Value = -12296.429
$ws.Range("H134").Value = 1798.04
$ws.Range("I134").Value = 1368.7222
$ws.Range("J134").Value = 2902
$ws.Range("K134").Value = 4106.1666
$ws.Range("L134").Value = 8706
$ws.Range("M134").Value = -1571.1666
$ws.Range("N134").Value = -13776

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 842.6667
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H87").Value = 7000
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 7000
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()
$ws.Range("M87").Value = 21000
$ws.Range("N87").Value = -23496
$ws.Range("H90").Value = 7000
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 7000
$ws.Range("K90").Value = 0
$ws.Range("L90").ClearContents()
$ws.Range("M90").Value = 63000
$ws.Range("N90").Value = -75480
$ws.Range("H107").Value = 649.125
$ws.Range("I107").Value = 672.5
$ws.Range("J107").Value = 641.3333
$ws.Range("K107").Value = 2017.5
$ws.Range("L107").Value = 1923.9999
$ws.Range("M107").Value = -97.5
$ws.Range("N107").Value = -5763.9999
$ws.Range("H135").Value = 842.6667
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 4317.8945
$ws.Range("I46").Value = 1041
$ws.Range("J46").Value = 4499.9443
$ws.Range("K46").Value = 1041
$ws.Range("L46").Value = 4499.9443
$ws.Range("M46").Value = -885
$ws.Range("N46").Value = -4811.9443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1772.1666
$ws.Range("I7").Value = 1261.0769
$ws.Range("J7").Value = 3101
$ws.Range("K7").Value = 1261.0769
$ws.Range("L7").Value = 3101
$ws.Range("M7").Value = -1149.0769
$ws.Range("N7").Value = -3325
$ws.Range("H47").Value = 1557
$ws.Range("I47").Value = 1557
$ws.Range("K47").Value = 1557
$ws.Range("M47").Value = -1067
$ws.Range("H52").Value = 1557
$ws.Range("I52").Value = 1557
$ws.Range("K52").Value = 1557
$ws.Range("M52").Value = -1324
$ws.Range("H61").Value = 2538.6667
$ws.Range("I61").Value = 2161
$ws.Range("K61").Value = 2161
$ws.Range("M61").Value = -1959
$ws.Range("H68").Value = 2361.963
$ws.Range("I68").Value = 1936.7826
$ws.Range("J68").Value = 2677.4194
$ws.Range("K68").Value = 1936.7826
$ws.Range("L68").Value = 2677.4194
$ws.Range("M68").Value = -1187.7826
$ws.Range("N68").Value = -4175.419400000001
$ws.Range("H71").Value = 2361.963
$ws.Range("I71").Value = 1936.7826
$ws.Range("J71").Value = 2677.4194
$ws.Range("K71").Value = 9683.913
$ws.Range("L71").Value = 13387.097
$ws.Range("M71").Value = -5939.913
$ws.Range("N71").Value = -20875.097
$ws.Range("H100").Value = 1810.7222
$ws.Range("I100").Value = 1392
$ws.Range("J100").Value = 2648.1667
$ws.Range("K100").Value = 1392
$ws.Range("L100").Value = 2648.1667
$ws.Range("M100").Value = -851
$ws.Range("N100").Value = -3730.1667
$ws.Range("H113").Value = 2538.6667
$ws.Range("I113").Value = 2161
$ws.Range("K113").Value = 2161
$ws.Range("M113").Value = 9
$ws.Range("H126").Value = 1772.1666
$ws.Range("I126").Value = 1261.0769
$ws.Range("J126").Value = 3101
$ws.Range("K126").Value = 3783.2307
$ws.Range("L126").Value = 9303
$ws.Range("M126").Value = -1313.2307
$ws.Range("N126").Value = -14243

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2169.3333
$ws.Range("I100").Value = 2170.6667
$ws.Range("J100").Value = 2166.6667
$ws.Range("K100").Value = 4341.3334
$ws.Range("L100").Value = 4333.3334
$ws.Range("M100").Value = -3800.3334
$ws.Range("N100").Value = -5415.3334
$ws.Range("H126").Value = 1421.6
$ws.Range("I126").Value = 1221.8462
$ws.Range("J126").Value = 1792.5714
$ws.Range("K126").Value = 3665.5386
$ws.Range("L126").Value = 5377.7142
$ws.Range("M126").Value = -1195.5386
$ws.Range("N126").Value = -10317.7142
$ws.Range("H136").Value = 2507.756
$ws.Range("I136").Value = 1831.7587
$ws.Range("J136").Value = 4141.4165
$ws.Range("K136").Value = 5495.2761
$ws.Range("L136").Value = 12424.2495
$ws.Range("M136").Value = -2945.2761
$ws.Range("N136").Value = -17524.2495
